# The template presentation shipped with one placeholder slide (a
# title/subtitle layout) that was only there as starter content. The
# PPT-generation code now builds every slide (including images) itself,
# so the template should start out empty, with zero slides.
$p = $ppt.ActivePresentation

# Remove every slide from the deck (there is exactly one in this
# template), working backwards so indices stay valid.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $p.Slides.Item($i).Delete()
}
